# Update league base data (Australia ALeague) - 17-06-2024 21:10
# The source rows for three pairs of matches had been mixed up; this swaps
# back the full record (columns B..AD) between each pair of rows while
# leaving column A (the sequential row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("B$rowA`:AD$rowA")
    $rangeB = $sheet.Range("B$rowB`:AD$rowB")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-Rows $ws 73 74
Swap-Rows $ws 104 105
Swap-Rows $ws 124 125
